$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2831
$ws.Range("F3").Value = 1145
$ws.Range("F4").Value = 20714
$ws.Range("F6").Value = 2673
$ws.Range("F7").Value = 789
$ws.Range("F8").Value = 615
$ws.Range("F9").Value = 496
$ws.Range("F10").Value = 740
$ws.Range("F11").Value = 273
$ws.Range("F12").Value = 261
$ws.Range("F13").Value = 72
$ws.Range("F14").Value = 103
$ws.Range("F15").Value = 503
$ws.Range("F17").Value = 248
$ws.Range("F18").Value = 8
$ws.Range("F19").Value = 411
$ws.Range("F20").Value = 22
$ws.Range("F23").Value = 115

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 16
$ws.Range("F3").Value = 27
$ws.Range("F5").Value = 321
$ws.Range("F8").Value = 17
$ws.Range("F9").Value = 12
$ws.Range("F10").Value = 13
$ws.Range("F11").Value = 2
$ws.Range("F12").Value = 95
$ws.Range("F13").Value = 40
$ws.Range("F14").Value = 131
$ws.Range("F18").Value = 2
$ws.Range("F20").Value = 37
$ws.Range("F21").Value = 23
$ws.Range("F22").Value = 38

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6101
$ws.Range("F3").Value = 687
$ws.Range("F4").Value = 664
$ws.Range("F5").Value = 1484
$ws.Range("F6").Value = 47

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6101
$ws.Range("F4").Value = 664
$ws.Range("F5").Value = 1484
$ws.Range("F6").Value = 2831
$ws.Range("F7").Value = 1145
$ws.Range("F8").Value = 20714
$ws.Range("F9").Value = 16
$ws.Range("F10").Value = 27
$ws.Range("F11").Value = 96
$ws.Range("F12").Value = 123
$ws.Range("F13").Value = 321
$ws.Range("F14").Value = 2673
$ws.Range("F15").Value = 789
$ws.Range("F17").Value = 47
$ws.Range("F19").Value = 496
$ws.Range("F20").Value = 740
$ws.Range("F21").Value = 273
$ws.Range("F22").Value = 261
$ws.Range("F23").Value = 9
$ws.Range("F24").Value = 72
$ws.Range("F25").Value = 17
$ws.Range("F26").Value = 12
$ws.Range("F27").Value = 103
$ws.Range("F28").Value = 13
$ws.Range("F29").Value = 2
$ws.Range("F30").Value = 503
$ws.Range("F31").Value = 95
$ws.Range("F32").Value = 179
$ws.Range("F34").Value = 248
$ws.Range("F35").Value = 131
$ws.Range("F37").Value = 8
$ws.Range("F38").Value = 411
$ws.Range("F39").Value = 4
$ws.Range("F40").Value = 22
$ws.Range("F42").Value = 2
$ws.Range("F43").Value = 5
$ws.Range("F45").Value = 2
$ws.Range("F46").Value = 20
$ws.Range("F48").Value = 23
$ws.Range("F49").Value = 38
$ws.Range("F50").Value = 115
